$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2627102.17
$ws.Range("C7").Value = -40.87202387519579
$ws.Range("D7").Value = 2646
$ws.Range("E7").Value = 2646
$ws.Range("F7").Value = 992.8579629629629
$ws.Range("G7").Value = 5.83147956427541
